# Update the "Förändrad" (Changed) date column (C) from 45178 to 45179
# for every data row (rows 2 through 505) on the active worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2:C505").Value = 45179
